$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2957.111
$ws.Range("I98").Value = 2451.75
$ws.Range("K98").Value = 2451.75
$ws.Range("M98").Value = -953.75

$ws.Range("H107").Value = 837.0833
$ws.Range("I107").Value = 822.2727
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 822.2727
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1097.7273
$ws.Range("N107").Value = -4840

$ws.Range("H122").Value = 2957.111
$ws.Range("I122").Value = 2451.75
$ws.Range("K122").Value = 7355.25
$ws.Range("M122").Value = -4905.25

$ws.Range("H137").Value = 558086.4399999999
$ws.Range("I137").Value = 3147.9092
$ws.Range("J137").Value = 879366.6
$ws.Range("K137").Value = 9443.7276
$ws.Range("L137").Value = 2638099.8
$ws.Range("M137").Value = -6893.7276
$ws.Range("N137").Value = -2643199.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13882.035
$ws.Range("I32").Value = 14444.487
$ws.Range("K32").Value = 14444.487
$ws.Range("M32").Value = -14157.487

$ws.Range("H34").Value = 37514
$ws.Range("I34").Value = 5000
$ws.Range("K34").Value = 5000
$ws.Range("M34").Value = -4729

$ws.Range("H61").Value = 6088.9585
$ws.Range("I61").Value = 2303.5833
$ws.Range("K61").Value = 2303.5833
$ws.Range("M61").Value = -2091.5833

$ws.Range("H74").Value = 4447.769
$ws.Range("I74").Value = 1666.4814
$ws.Range("J74").Value = 10705.667
$ws.Range("K74").Value = 1666.4814
$ws.Range("L74").Value = 10705.667
$ws.Range("M74").Value = -792.4813999999999
$ws.Range("N74").Value = -12453.667

$ws.Range("H77").Value = 4447.769
$ws.Range("I77").Value = 1666.4814
$ws.Range("J77").Value = 10705.667
$ws.Range("K77").Value = 8332.406999999999
$ws.Range("L77").Value = 53528.335
$ws.Range("M77").Value = -3964.406999999999
$ws.Range("N77").Value = -62264.335

$ws.Range("H122").Value = 1893.6818
$ws.Range("I122").Value = 1850.5238
$ws.Range("K122").Value = 5551.5714
$ws.Range("M122").Value = -3101.5714

$ws.Range("H132").Value = 4187.579
$ws.Range("I132").Value = 3644.3333
$ws.Range("J132").Value = 4438.3076
$ws.Range("K132").Value = 10932.9999
$ws.Range("L132").Value = 13314.9228
$ws.Range("M132").Value = -8402.999899999999
$ws.Range("N132").Value = -18374.9228

$ws.Range("H136").Value = 6088.9585
$ws.Range("I136").Value = 2303.5833
$ws.Range("K136").Value = 6910.749899999999
$ws.Range("M136").Value = -4360.749899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 58167.723
$ws.Range("I134").Value = 3135.9285
$ws.Range("J134").Value = 250779
$ws.Range("K134").Value = 9407.7855
$ws.Range("L134").Value = 752337
$ws.Range("M134").Value = -6872.7855
$ws.Range("N134").Value = -757407

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 480563.97
$ws.Range("I31").Value = 8902.166999999999
$ws.Range("J31").Value = 702522.4399999999
$ws.Range("K31").Value = 8902.166999999999
$ws.Range("L31").Value = 702522.4399999999
$ws.Range("M31").Value = -8607.166999999999
$ws.Range("N31").Value = -703112.4399999999

$ws.Range("H34").Value = 480563.97
$ws.Range("I34").Value = 8902.166999999999
$ws.Range("J34").Value = 702522.4399999999
$ws.Range("K34").Value = 8902.166999999999
$ws.Range("L34").Value = 702522.4399999999
$ws.Range("M34").Value = -8700.166999999999
$ws.Range("N34").Value = -702926.4399999999

$ws.Range("H132").Value = 3138.75
$ws.Range("I132").Value = 2869.4634
$ws.Range("J132").Value = 4716
$ws.Range("K132").Value = 8608.3902
$ws.Range("L132").Value = 14148
$ws.Range("M132").Value = -6078.3902
$ws.Range("N132").Value = -19208

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 248
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 248
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 744
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -1082

$ws.Range("H34").Value = 2531.476
$ws.Range("J34").Value = 3053.1177
$ws.Range("L34").Value = 9159.3531
$ws.Range("N34").Value = -9327.3531

$ws.Range("H39").Value = 5766.6665
$ws.Range("J39").Value = 5766.6665
$ws.Range("L39").Value = 17299.9995
$ws.Range("N39").Value = -17887.9995

$ws.Range("H55").Value = 3700
$ws.Range("I55").Value = 1900
$ws.Range("J55").Value = 4300
$ws.Range("K55").Value = 5700
$ws.Range("L55").Value = 12900
$ws.Range("M55").Value = -5523
$ws.Range("N55").Value = -13254

$ws.Range("H68").Value = 1649.2142
$ws.Range("I68").Value = 1481.4043
$ws.Range("J68").Value = 1803.8628
$ws.Range("K68").Value = 4444.2129
$ws.Range("L68").Value = 5411.588400000001
$ws.Range("M68").Value = -3633.2129
$ws.Range("N68").Value = -7033.588400000001

$ws.Range("H71").Value = 1649.2142
$ws.Range("I71").Value = 1481.4043
$ws.Range("J71").Value = 1803.8628
$ws.Range("K71").Value = 13332.6387
$ws.Range("L71").Value = 16234.7652
$ws.Range("M71").Value = -9276.6387
$ws.Range("N71").Value = -24346.7652

$ws.Range("H137").Value = 28108.7
$ws.Range("I137").Value = 1720.5333
$ws.Range("J137").Value = 107273.2
$ws.Range("K137").Value = 5161.5999
$ws.Range("L137").Value = 321819.6
$ws.Range("M137").Value = -61.59990000000016
$ws.Range("N137").Value = -332019.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2141.6667
$ws.Range("I113").Value = 2063.6365
$ws.Range("K113").Value = 2063.6365
$ws.Range("M113").Value = 106.3634999999999

$ws.Range("H132").Value = 15193.6
$ws.Range("I132").Value = 8012
$ws.Range("J132").Value = 16989
$ws.Range("K132").Value = 24036
$ws.Range("L132").Value = 50967
$ws.Range("M132").Value = -21506
$ws.Range("N132").Value = -56027

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 36666.668
$ws.Range("J25").Value = 36666.668
$ws.Range("L25").Value = 36666.668
$ws.Range("N25").Value = -37126.668

$ws.Range("H111").Value = 48000
$ws.Range("J111").Value = 48000
$ws.Range("L111").Value = 48000
$ws.Range("N111").Value = -56180

$ws.Range("H132").Value = 7090.811
$ws.Range("I132").Value = 9398.208000000001
$ws.Range("J132").Value = 2831
$ws.Range("K132").Value = 28194.624
$ws.Range("L132").Value = 8493
$ws.Range("M132").Value = -25664.624
$ws.Range("N132").Value = -13553

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 30000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 30000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 30000
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -30580

$ws.Range("H40").Value = 59000
$ws.Range("J40").Value = 59000
$ws.Range("L40").Value = 59000
$ws.Range("N40").Value = -59298
